$wb = $excel.ActiveWorkbook

# The sheet that needs editing is "Aerobox" (the active/selected sheet in the file)
$ws = $wb.Worksheets.Item("Aerobox")
$ws.Activate()

# Update the CUPO (quota) values
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 8

# Update the selection to match the saved cursor position
$ws.Range("C2").Select()
